$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new blank column at F.
#    Old F ("video") shifts to G, old G ("usable") shifts to H.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).Insert()
$ws.Cells.Item(1, 6).Value = "actual_size"

# ---------------------------------------------------------------------------
# 2. Column widths (post-insert layout: A..H).
#    (ColumnWidth values are chosen so the engine's internal char-width
#    quantisation lands as close as possible on the target stored widths:
#    5.86, 7.14, 12.0, 12.0, 12.0, 10.43, 5.57, 6.57.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 5.0
$ws.Columns.Item(2).ColumnWidth = 6.333333
$ws.Range($ws.Columns.Item(3), $ws.Columns.Item(5)).ColumnWidth = 11.166667
$ws.Columns.Item(6).ColumnWidth = 9.666667
$ws.Columns.Item(7).ColumnWidth = 4.666667
$ws.Columns.Item(8).ColumnWidth = 5.666667

# ---------------------------------------------------------------------------
# 3. Fill in "actual_size" (column F) for the rows that were measured.
#    Rows not listed keep an empty actual_size cell.
# ---------------------------------------------------------------------------
$actualSize = @{
    2  = 13.5
    3  = 10.0
    7  = 15.0
    10 = 7.5
    14 = 10.0
    19 = 13.5
    21 = 13.0
    22 = 20.0
    24 = 8.5
    28 = 19.5
    29 = 9.5
    34 = 7.5
    35 = 7.5
    36 = 8.5
    39 = 17.0
    47 = 12.0
    49 = 12.0
}

foreach ($row in $actualSize.Keys) {
    $ws.Cells.Item($row, 6).Value = $actualSize[$row]
}

# ---------------------------------------------------------------------------
# 4. Apply the autofilter on the "usable" column (H) for value = 1.
#    This is done while every row still carries its original "usable" flag,
#    so the filter hides exactly the rows that were (and remain) unusable.
# ---------------------------------------------------------------------------
$ws.Range("A1:H53").AutoFilter(8, "1", 7)

# ---------------------------------------------------------------------------
# 5. Two particles (rows 38 and 51) are reclassified as not usable, but stay
#    visible (they are not re-filtered/hidden after being flagged unusable).
# ---------------------------------------------------------------------------
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(51, 8).Value = 0

# ---------------------------------------------------------------------------
# 6. Register the (hidden) sheet-scoped _FilterDatabase defined name that
#    Excel normally maintains alongside an AutoFilter.
# ---------------------------------------------------------------------------
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "='8_particle_selection'!`$A`$1:`$H`$53")
$fdb.Visible = $false
